{"js": "// Ordered list of (find, replace) pairs mirroring the source diff: the\n// worksheet date banner plus every two-digit x two-digit multiplication\n// prompt in the practice table. Every \"find\" string is unique within the\n// document, so searching the whole body for each one and swapping its text\n// in place is unambiguous and leaves formatting (fonts/sizes), the\n// paragraph/table structure, and every other cell untouched.\nconst pairs = [\n    [\"2025-05-22 Thursday\", \"2025-05-23 Friday\"],\n    [\"56\u00d754=\", \"80\u00d729=\"],\n    [\"32\u00d771=\", \"25\u00d722=\"],\n    [\"43\u00d723=\", \"89\u00d797=\"],\n    [\"84\u00d752=\", \"41\u00d743=\"],\n    [\"99\u00d760=\", \"88\u00d794=\"],\n    [\"82\u00d796=\", \"12\u00d749=\"],\n    [\"67\u00d726=\", \"45\u00d766=\"],\n    [\"17\u00d715=\", \"43\u00d757=\"],\n    [\"13\u00d713=\", \"39\u00d752=\"],\n    [\"86\u00d726=\", \"98\u00d771=\"],\n    [\"93\u00d729=\", \"40\u00d735=\"],\n    [\"26\u00d784=\", \"69\u00d768=\"],\n    [\"40\u00d787=\", \"65\u00d798=\"],\n    [\"43\u00d730=\", \"12\u00d742=\"],\n    [\"50\u00d753=\", \"11\u00d730=\"],\n    [\"64\u00d741=\", \"40\u00d750=\"],\n    [\"74\u00d758=\", \"20\u00d747=\"],\n    [\"17\u00d735=\", \"50\u00d742=\"],\n    [\"97\u00d725=\", \"65\u00d784=\"],\n    [\"69\u00d726=\", \"77\u00d714=\"],\n    [\"67\u00d780=\", \"24\u00d747=\"],\n    [\"69\u00d780=\", \"75\u00d763=\"],\n    [\"26\u00d743=\", \"56\u00d749=\"],\n    [\"36\u00d757=\", \"56\u00d717=\"],\n    [\"16\u00d722=\", \"36\u00d730=\"],\n];\n\nfor (const [findText, replaceText] of pairs) {\n    const results = context.document.body.search(findText, { matchCase: true });\n    results.load(\"text\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Ordered list of (find, replace) pairs mirroring the source diff: the\n# worksheet date banner plus every two-digit x two-digit multiplication\n# prompt in the practice table. Each \"find\" string is unique in the\n# document, so a simple whole-document Find/Replace for each pair is\n# unambiguous and leaves every other run (fonts, sizes, paragraph/table\n# structure) untouched.\n$pairs = @(\n    @(\"2025-05-22 Thursday\", \"2025-05-23 Friday\"),\n    @(\"56\u00d754=\", \"80\u00d729=\"),\n    @(\"32\u00d771=\", \"25\u00d722=\"),\n    @(\"43\u00d723=\", \"89\u00d797=\"),\n    @(\"84\u00d752=\", \"41\u00d743=\"),\n    @(\"99\u00d760=\", \"88\u00d794=\"),\n    @(\"82\u00d796=\", \"12\u00d749=\"),\n    @(\"67\u00d726=\", \"45\u00d766=\"),\n    @(\"17\u00d715=\", \"43\u00d757=\"),\n    @(\"13\u00d713=\", \"39\u00d752=\"),\n    @(\"86\u00d726=\", \"98\u00d771=\"),\n    @(\"93\u00d729=\", \"40\u00d735=\"),\n    @(\"26\u00d784=\", \"69\u00d768=\"),\n    @(\"40\u00d787=\", \"65\u00d798=\"),\n    @(\"43\u00d730=\", \"12\u00d742=\"),\n    @(\"50\u00d753=\", \"11\u00d730=\"),\n    @(\"64\u00d741=\", \"40\u00d750=\"),\n    @(\"74\u00d758=\", \"20\u00d747=\"),\n    @(\"17\u00d735=\", \"50\u00d742=\"),\n    @(\"97\u00d725=\", \"65\u00d784=\"),\n    @(\"69\u00d726=\", \"77\u00d714=\"),\n    @(\"67\u00d780=\", \"24\u00d747=\"),\n    @(\"69\u00d780=\", \"75\u00d763=\"),\n    @(\"26\u00d743=\", \"56\u00d749=\"),\n    @(\"36\u00d757=\", \"56\u00d717=\"),\n    @(\"16\u00d722=\", \"36\u00d730=\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    [void]$find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n}\n"}
